$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns (price / volume change) stay as plain text,
# matching the inline-string cell type used throughout this sheet,
# so Excel does not auto-convert values like "1.00" or "12.27" into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "69.022.29"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "3.742.84"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "601.60"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "168.24"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").Value = "3.741.31"
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("E11").Value = "  +2.87%  "
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("D13").Value = "37.95"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "4.365.09"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "3.741.22"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "69.032.17"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "7.28"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("D20").Value = "17.19"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").Value = "10.75"
$ws.Range("E21").Value = "  +16.30%  "
$ws.Range("D22").Value = "492.44"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("E24").Value = "  +5.49%  "
$ws.Range("D25").Value = "84.72"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").Value = "12.27"
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  -0.31%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "2.97"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("E31").Value = "  +5.06%  "
$ws.Range("D32").Value = "8.01"
$ws.Range("E32").Value = "  +3.09%  "
$ws.Range("D33").Value = "31.48"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "3.884.84"
$ws.Range("E34").Value = "  +0.51%  "
$ws.Range("D35").Value = "0.109"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").Value = "3.672.79"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("D39").Value = "5.84"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").Value = "2.95"
$ws.Range("E42").Value = "  +3.21%  "
$ws.Range("D43").Value = "429.60"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "48.54"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "8.47"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D50").Value = "2.780.71"
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("E51").Value = "  +0.04%  "

# Rows 48 and 49 swap content: Arweave <-> Monero positions
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "141.63"
$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "39.99"
$ws.Range("E49").Value = "  -1.82%  "
